$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new product row (row 15)
$ws.Range("A15").Value = "can5"
$ws.Range("B15").Value = "canet"
$ws.Range("C15").Value = 1200

# Clear the style previously applied to column A (rows 2-14), leaving the
# header row (row 1) untouched, so those cells revert to the default style.
$ws.Range("A2:A14").Style = "Normal"
